# Reproduce the author's edit described in the commit ("commit from my work
# computer"): on Sheet1 (the "表2" trading log table), the 卖出价 (sell price)
# recorded for the 2021-02-18 row (row 6) was corrected from 93.3 to 85.
# Every other cell that differs in the diff (K6/L6/M6/N6, and B/K/N for all
# the following blank rows) is a table-calculated column that is simply
# recomputed downstream of this one input cell - they are not separate edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H6 = 卖出价 (sell price) for the row dated 20210218
$ws.Range("H6").Value = 85

# The author's cursor ended up on I8 (next day's 卖出数量 entry) when the
# workbook was saved from the other machine.
$ws.Range("I8").Select()
